$wb = $excel.ActiveWorkbook

# --- Work on the "Chests" sheet: add new rows for "Your Cave" ---
$chests = $wb.Worksheets.Item("Chests")

# Fill rows in the same order the shared strings were first introduced
# (134, 135, 133, 136, 137) so the resulting shared-string table order
# matches the target workbook exactly.
$chests.Range("A4").Value = 134
$chests.Range("B4").Value = "Your Cave (459)"
$chests.Range("C4").Value = "1x Levitation, 1x Healing Potion II"

$chests.Range("A5").Value = 135
$chests.Range("B5").Value = "Your Cave (460)"
$chests.Range("C5").Value = "3x Torch"

$chests.Range("A3").Value = 133
$chests.Range("B3").Value = "Your Cave (459)"
$chests.Range("C3").Value = "2x Rope"

$chests.Range("A6").Value = 136
$chests.Range("B6").Value = "Your Cave (459)"
$chests.Range("C6").Value = "1x Rope"

$chests.Range("A7").Value = 137
$chests.Range("B7").Value = "Your Cave (459)"
$chests.Range("C7").Value = "1x Rope"

# Select C8 on the Chests sheet, and make it the active sheet/selection
$chests.Activate()
$chests.Range("C8").Select()

$wb.Save()
